$d = $word.ActiveDocument

# The document ends with a paragraph whose last run reads "Fixed! Very simple".
# The commit adds, right after that paragraph:
#   1. an empty paragraph (same tab-stop/spacing formatting)
#   2. a paragraph with the text "V0.9.4.3 Beta"
#   3. a paragraph with the text "Added in Game Page Layout as default config in backend"

$pFixed = $d.Paragraphs.Last

# 1) Blank paragraph right after "Fixed! Very simple".
#    InsertParagraphAfter() duplicates pFixed's paragraph formatting (tabs + spacing)
#    onto the new paragraph, which is exactly what the diff shows.
$pFixed.Range.InsertParagraphAfter()
$pBlank = $d.Paragraphs.Last

# 2) "V0.9.4.3 Beta" paragraph.
$pBlank.Range.InsertParagraphAfter()
$pBeta = $d.Paragraphs.Last
$pBeta.Range.Text = "V0.9.4.3 Beta"

# 3) "Added in Game Page Layout as default config in backend" paragraph.
$pBeta.Range.InsertParagraphAfter()
$pAdded = $d.Paragraphs.Last
$pAdded.Range.Text = "Added in Game Page Layout as default config in backend"

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
